$wb = $excel.ActiveWorkbook

$wsOffline = $wb.Worksheets.Item("Thiết bị Offline")
$wsOnline  = $wb.Worksheets.Item("Thiết bị Online")

# --- Update header row (row 9) and template row (row 10) for both sheets ---
foreach ($ws in @($wsOffline, $wsOnline)) {
    $ws.Range("H9").Value = "IpAddress"
    $ws.Range("G10").Value = "{{item.CheckConnectTime}}"
    $ws.Range("H10").Value = "{{item.IPAddress}}"

    # Remove column I entirely (shift cells left)
    $ws.Range("I9:I10").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
}

# --- Defined Names ---
$wb.Names.Item("datas").RefersTo = "='Thiết bị Online'!`$A`$10:`$H`$11"
$wb.Names.Item("items").RefersTo = "='Thiết bị Offline'!`$A`$10:`$H`$11"

$wsOffline.PageSetup.PrintArea = "`$A`$1:`$N`$22"
$wsOnline.PageSetup.PrintArea = "`$A`$1:`$I`$21"

# --- Page setup scale change on offline sheet ---
$wsOffline.PageSetup.Zoom = $false
$wsOffline.PageSetup.FitToPagesWide = $false
$wsOffline.PageSetup.FitToPagesTall = $false
$wsOffline.PageSetup.Zoom = 57

# --- Selection changes (cosmetic) ---
$wsOnline.Activate() | Out-Null
$wsOnline.Range("H25").Select() | Out-Null
$wsOffline.Activate() | Out-Null
$wsOffline.Range("E31").Select() | Out-Null
